$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$powerOffText = "Button (LV_EVENT_LONG_PRESSED) EVT_BUTTON_BACK_LONG_PRESSED"

# Insert three new rows (for the new "power off" long-press action) after
# rows 3, 4 and 5 (pre-insert numbering), i.e. right after Browser_Disk,
# Browser_NoDisk and Browser_BadDisk respectively.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(8).Insert()

# Row 4: new long-press power off row (no State name)
$ws.Cells.Item(4, 1).Clear()
$ws.Cells.Item(4, 3).Value = $powerOffText

# Row 6: new long-press power off row
$ws.Cells.Item(6, 1).Clear()
$ws.Cells.Item(6, 3).Value = $powerOffText

# Row 8: new long-press power off row
$ws.Cells.Item(8, 1).Clear()
$ws.Cells.Item(8, 3).Value = $powerOffText

# Resize the table to include the new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:F13"))

$ws.Range("G7").Select()
